# "Loan RBI, Variable Instalments"
#
# The "Repayment schedule" sheet gets a new, blank column inserted in front
# of what used to be column N ("Late"); that pushes the old N/O/P data
# (Late / heading / Outstanding) one slot to the right, to O/P/Q. The
# "Repayment schedule" tab also becomes the active sheet/tab (it was
# "Transactions" before), with the cell cursor left on J19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before (old) column N - this shifts
# N:P -> O:Q, copying each cell's style from the column to the left,
# exactly like Excel's native "Insert Column" command.
[void]$ws.Columns("N").Insert()

# The freshly inserted column picks up column M's width in real Excel.
# Re-apply that width (closest representable value).
$ws.Columns("N").ColumnWidth = 9.8

# "Repayment schedule" becomes the active sheet (was "Transactions"),
# and the selection on it moves to J19.
[void]$ws.Activate()
[void]$ws.Range("J19").Select()
